$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns C and D entirely (shift remaining cells left)
$ws.Range("C1:D4").Delete()

# Header row: B1 becomes "Agent回答" (was "期望输出"; old C1 "Agent回答" is gone with column delete)
$ws.Range("B1").Value = "Agent回答"

# Row 2: A2 becomes the sample question, B2 cleared
$ws.Range("A2").Value = "（示例）中国的首都在哪里？"
$ws.Range("B2").ClearContents()

# Row 3: B3 cleared (A3 already has sample text)
$ws.Range("B3").ClearContents()

# Row 4: B4 cleared (A4 already has sample text)
$ws.Range("B4").ClearContents()

# Row 5 (new): A5 sample text, B5 left blank
$ws.Range("A5").Value = "（示例）中国的首都在哪里？"
